$d = $word.ActiveDocument

# --- Rename "Chance Morley" -> "Highmage Tarius" throughout the document ---

# Body-paragraph mention: " Chance Morley is currently..." -> " Highmage Tarius currently..."
$d.Content.Find.Execute('Chance Morley is', $true, $false, $false, $false, $false, $true, 1, $false, 'Highmage Tarius', 2)

# Heading mention: "Chance Morley" -> "Highmage Tarius"
$d.Content.Find.Execute('Chance Morley', $true, $false, $false, $false, $false, $true, 1, $false, 'Highmage Tarius', 2)

# "None has really seen Chance use" -> "...Tarius use"
$d.Content.Find.Execute('Chance use', $true, $false, $false, $false, $false, $true, 1, $false, 'Tarius use', 2)

# "...that it was Chance work" -> "...that it was Tarius work"
$d.Content.Find.Execute('Chance work', $true, $false, $false, $false, $false, $true, 1, $false, 'Tarius work', 2)

# --- Re-join a handful of unrelated mid-word run splits elsewhere in the document ---
# (no textual change -- replace each anchor with itself so Word coalesces the runs
#  it spans, matching the canonical save produced alongside the rename above)
$d.Content.Find.Execute(' them fresh water. That’s why the humans', $true, $false, $false, $false, $false, $true, 1, $false, ' them fresh water. That’s why the humans', 2)
$d.Content.Find.Execute(' was captured by the Trol', $true, $false, $false, $false, $false, $true, 1, $false, ' was captured by the Trol', 2)
$d.Content.Find.Execute('re they left, now train other humans wit', $true, $false, $false, $false, $false, $true, 1, $false, 're they left, now train other humans wit', 2)
$d.Content.Find.Execute('ing the invasion of Tendora, Gavin Burbr', $true, $false, $false, $false, $false, $true, 1, $false, 'ing the invasion of Tendora, Gavin Burbr', 2)
$d.Content.Find.Execute('last fight against the trolls. Defeating', $true, $false, $false, $false, $false, $true, 1, $false, 'last fight against the trolls. Defeating', 2)
$d.Content.Find.Execute('rew up to become a soldi', $true, $false, $false, $false, $false, $true, 1, $false, 'rew up to become a soldi', 2)
$d.Content.Find.Execute('as the first knight. During the later pa', $true, $false, $false, $false, $false, $true, 1, $false, 'as the first knight. During the later pa', 2)
$d.Content.Find.Execute(' giving them faith that they could defea', $true, $false, $false, $false, $false, $true, 1, $false, ' giving them faith that they could defea', 2)

Write-Output "done"
